$d = $word.ActiveDocument

# The "Exerc. 1" paragraph originally contains leftover spell-check markup
# (<w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>)
# splitting "Exerc" and ". 1" into two runs. The target edit collapses this
# into a single clean run/paragraph and adds two new explanatory paragraphs
# right after it.

# Step 1: Insert a brand-new, attribute-free paragraph right after the blank
# paragraph that precedes "Exerc. 1", and give it the plain text "Exerc. 1".
$blankBefore = $d.Paragraphs(2)
$blankBefore.Range.InsertParagraphAfter()
$newExerc = $d.Paragraphs(3)
$newExerc.Range.Text = "Exerc. 1"

# Step 2: Remove the old "Exerc. 1" paragraph (now pushed to index 4), which
# still carries the stray proofErr markers and split runs, paragraph mark
# included, so it disappears completely.
$oldExerc = $d.Paragraphs(4)
$oldExerc.Range.Delete()

# Step 3: Insert the two new explanatory paragraphs right after the cleaned
# up "Exerc. 1" paragraph (index 3).
$exercPara = $d.Paragraphs(3)
$exercPara.Range.InsertParagraphAfter()
$note1 = $d.Paragraphs(4)
$note1.Range.Text = "Call the attention that it is Active Directory under Settings for the Synapse workspace. Many students might go straight away to the first Active Directory link they see in the left panel."

$note1.Range.InsertParagraphAfter()
$note2 = $d.Paragraphs(5)
$note2.Range.Text = "It is System Managed Identity, not User."
